$d = $word.ActiveDocument

# 1. Remove the "Meta description" paragraph that immediately follows the
#    opening "Play Candyfinity Online Slot Game for Free!" heading.
$metaPara = $d.Paragraphs.Item(2)
$metaPara.Range.Delete()

# 2. Insert a new bold "Play Candyfinity Online Slot Game for Free!" paragraph
#    right before the closing paragraph (the one with the italic image-prompt
#    text). InsertXML merges the trailing paragraph of the inserted fragment
#    with the target paragraph, so we include an extra empty paragraph in the
#    fragment and then delete the stray empty paragraph it leaves behind.
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$insertAt = $d.Range($lastPara.Range.Start, $lastPara.Range.Start)
$insertAt.InsertXML("<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Candyfinity Online Slot Game for Free!</w:t></w:r></w:p><w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'/>")

$strayPara = $d.Paragraphs.Item($lastPara.Index + 1)
$strayPara.Range.Delete()

# 3. Swap the old image-generation prompt text for the new meta-description
#    text inside the final (italic) run, keeping its run formatting intact.
$d.Content.Find.Execute("Create an eye-catching feature image for Candyfinity in cartoon style featuring a happy Maya warrior with glasses. The warrior should be surrounded by various gummy candies, lollipops, sugary glazes, and all kinds of sweets that are the ingredients for this tasty and exciting slot game. Use vibrant colors such as red, pink, and violet to capture the explosion of colors and shapes in the game. The image should also include the game logo " + [char]34 + "Candyfinity" + [char]34 + " in bold and playful font. Let your creativity shine to grab the attention of online slot game players looking for a fun and exciting game to play.", $true, $false, $false, $false, $false, $true, 1, $false, "Read our review of Candyfinity - a candy-themed slot game with various exciting features and winning potential. Play for free today!", 2)
